$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.01509998495488674 ; $ws.Range("D2").Value = 0.1667493158925453 ; $ws.Range("E2").Value = 0.2567807369960775 ; $ws.Range("F2").Value = 1.160573058580006 ; $ws.Range("G2").Value = 0.8432839204109825 ; $ws.Range("H2").Value = 0.8743153596641946 ; $ws.Range("I2").Value = 0.7479278897184827 ; $ws.Range("J2").Value = 0.4192787400100997 ; $ws.Range("M2").Value = 10.9622247095927
$ws.Range("C3").Value = 0.01453200457124382 ; $ws.Range("D3").Value = 0.1565091770063134 ; $ws.Range("E3").Value = 0.2336856282998809 ; $ws.Range("F3").Value = 1.212554290379188 ; $ws.Range("G3").Value = 0.8626374261524177 ; $ws.Range("H3").Value = 0.8988939228847528 ; $ws.Range("I3").Value = 0.7627294601680461 ; $ws.Range("J3").Value = 0.3762519443778274 ; $ws.Range("M3").Value = 9.645719339023458
$ws.Range("C4").Value = 0.01419861944958711 ; $ws.Range("D4").Value = 0.1502392157268986 ; $ws.Range("E4").Value = 0.2196314938467196 ; $ws.Range("F4").Value = 1.247343426654027 ; $ws.Range("G4").Value = 0.8768852520360468 ; $ws.Range("H4").Value = 0.9155250563511572 ; $ws.Range("I4").Value = 0.7735178989497911 ; $ws.Range("J4").Value = 0.35010256972555 ; $ws.Range("M4").Value = 8.834820259949765
$ws.Range("C5").Value = 0.0140665178752144 ; $ws.Range("D5").Value = 0.1476885769251908 ; $ws.Range("E5").Value = 0.2139346389647727 ; $ws.Range("F5").Value = 1.262228232202077 ; $ws.Range("G5").Value = 0.8832717464590161 ; $ws.Range("H5").Value = 0.9226831557936066 ; $ws.Range("I5").Value = 0.7783327691857451 ; $ws.Range("J5").Value = 0.339510221410535 ; $ws.Range("M5").Value = 8.503685594311264
$ws.Range("C6").Value = 0.01404480637508954 ; $ws.Range("D6").Value = 0.1472653133601369 ; $ws.Range("E6").Value = 0.2129904718546243 ; $ws.Range("F6").Value = 1.264742185929649 ; $ws.Range("G6").Value = 0.884366861250129 ; $ws.Range("H6").Value = 0.9238945649214116 ; $ws.Range("I6").Value = 0.7791572890568901 ; $ws.Range("J6").Value = 0.3377551090127326 ; $ws.Range("M6").Value = 8.448658298327871
$ws.Range("C7").Value = 0.01419682280996426 ; $ws.Range("D7").Value = 0.1502047989767732 ; $ws.Range("E7").Value = 0.2195545431317996 ; $ws.Range("F7").Value = 1.247541320664041 ; $ws.Range("G7").Value = 0.8769690511858954 ; $ws.Range("H7").Value = 0.9156200596305553 ; $ws.Range("I7").Value = 0.7735811513408919 ; $ws.Range("J7").Value = 0.3499594645880109 ; $ws.Range("M7").Value = 8.83035729583321
$ws.Range("C8").Value = 0.01490089162074071 ; $ws.Range("D8").Value = 0.1632148744426303 ; $ws.Range("E8").Value = 0.2487903938478411 ; $ws.Range("F8").Value = 1.177891311432866 ; $ws.Range("G8").Value = 0.849457582398955 ; $ws.Range("H8").Value = 0.8824665755510921 ; $ws.Range("I8").Value = 0.7526731091818561 ; $ws.Range("J8").Value = 0.4043847226512014 ; $ws.Range("M8").Value = 10.50879697590932
$ws.Range("C9").Value = 0.01640866313633893 ; $ws.Range("D9").Value = 0.1888690344609074 ; $ws.Range("E9").Value = 0.3071999425101239 ; $ws.Range("F9").Value = 1.064752511115032 ; $ws.Range("G9").Value = 0.8149209456678221 ; $ws.Range("H9").Value = 0.8299623102788445 ; $ws.Range("I9").Value = 0.7255759725060713 ; $ws.Range("J9").Value = 0.513448377316962 ; $ws.Range("M9").Value = 13.78226485224189
$ws.Range("C10").Value = 0.01760186035081546 ; $ws.Range("D10").Value = 0.2078111265351765 ; $ws.Range("E10").Value = 0.3508902711183026 ; $ws.Range("F10").Value = 0.9968375219392982 ; $ws.Range("G10").Value = 0.8022975816409144 ; $ws.Range("H10").Value = 0.7994230584910724 ; $ws.Range("I10").Value = 0.7147250382995907 ; $ws.Range("J10").Value = 0.5953084369811279 ; $ws.Range("M10").Value = 16.18022664409892
$ws.Range("C11").Value = 0.01816514429425808 ; $ws.Range("D11").Value = 0.2164513094042491 ; $ws.Range("E11").Value = 0.3709643270364467 ; $ws.Range("F11").Value = 0.9694482580049026 ; $ws.Range("G11").Value = 0.799528438829185 ; $ws.Range("H11").Value = 0.78736622679682 ; $ws.Range("I11").Value = 0.7118853405888785 ; $ws.Range("J11").Value = 0.6329988810644807 ; $ws.Range("M11").Value = 17.27056655318995
$ws.Range("C12").Value = 0.0183815612331415 ; $ws.Range("D12").Value = 0.2197266915462421 ; $ws.Range("E12").Value = 0.378597065916324 ; $ws.Range("F12").Value = 0.9595983540351654 ; $ws.Range("G12").Value = 0.7989248775939757 ; $ws.Range("H12").Value = 0.7830724104677529 ; $ws.Range("I12").Value = 0.711122570383651 ; $ws.Range("J12").Value = 0.6473427580752968 ; $ws.Range("M12").Value = 17.68346161158689
$ws.Range("C13").Value = 0.01833481072044663 ; $ws.Range("D13").Value = 0.2190211190474542 ; $ws.Range("E13").Value = 0.3769517926038759 ; $ws.Range("F13").Value = 0.9616962111442007 ; $ws.Range("G13").Value = 0.7990347869888694 ; $ws.Range("H13").Value = 0.7839849422094289 ; $ws.Range("I13").Value = 0.7112727650107047 ; $ws.Range("J13").Value = 0.6442502695490191 ; $ws.Range("M13").Value = 17.59453561829969
$ws.Range("C14").Value = 0.01818288565033299 ; $ws.Range("D14").Value = 0.2167207051147386 ; $ws.Range("E14").Value = 0.3715916381482884 ; $ws.Range("F14").Value = 0.9686273390829001 ; $ws.Range("G14").Value = 0.7994697618008502 ; $ws.Range("H14").Value = 0.7870074777910077 ; $ws.Range("I14").Value = 0.7118162580688647 ; $ws.Range("J14").Value = 0.6341774930824045 ; $ws.Range("M14").Value = 17.30453505634006
$ws.Range("C15").Value = 0.01809023777650509 ; $ws.Range("D15").Value = 0.215312101117263 ; $ws.Range("E15").Value = 0.3683125203380655 ; $ws.Range("F15").Value = 0.9729413461126768 ; $ws.Range("G15").Value = 0.7997946840604868 ; $ws.Range("H15").Value = 0.788894508259034 ; $ws.Range("I15").Value = 0.7121902030949201 ; $ws.Range("J15").Value = 0.6280171110579715 ; $ws.Range("M15").Value = 17.12690473414119
$ws.Range("C16").Value = 0.0175654747758216 ; $ws.Range("D16").Value = 0.2072469540159716 ; $ws.Range("E16").Value = 0.3495825998969053 ; $ws.Range("F16").Value = 0.9986995181955791 ; $ws.Range("G16").Value = 0.8025398653117861 ; $ws.Range("H16").Value = 0.8002486068158703 ; $ws.Range("I16").Value = 0.7149537458900284 ; $ws.Range("J16").Value = 0.5928549029896999 ; $ws.Range("M16").Value = 16.10896685032543
$ws.Range("C17").Value = 0.01724891335767609 ; $ws.Range("D17").Value = 0.2023053369495926 ; $ws.Range("E17").Value = 0.3381451110734446 ; $ws.Range("F17").Value = 1.015411623554556 ; $ws.Range("G17").Value = 0.804997757231007 ; $ws.Range("H17").Value = 0.8076896619134288 ; $ws.Range("I17").Value = 0.7171938294780063 ; $ws.Range("J17").Value = 0.571404139382139 ; $ws.Range("M17").Value = 15.48442002179331
$ws.Range("C18").Value = 0.01706875379779405 ; $ws.Range("D18").Value = 0.1994652324355002 ; $ws.Range("E18").Value = 0.3315851492935167 ; $ws.Range("F18").Value = 1.02535285708764 ; $ws.Range("G18").Value = 0.8066902536951659 ; $ws.Range("H18").Value = 0.8121418526945376 ; $ws.Range("I18").Value = 0.7186789105907181 ; $ws.Range("J18").Value = 0.5591083552816087 ; $ws.Range("M18").Value = 15.12514448472109
$ws.Range("C19").Value = 0.01700807942354743 ; $ws.Range("D19").Value = 0.1985039919918847 ; $ws.Range("E19").Value = 0.3293671770842508 ; $ws.Range("F19").Value = 1.028774741294022 ; $ws.Range("G19").Value = 0.8073106605572065 ; $ws.Range("H19").Value = 0.8136786430443976 ; $ws.Range("I19").Value = 0.719215174750083 ; $ws.Range("J19").Value = 0.5549522635855624 ; $ws.Range("M19").Value = 15.00348859574927
$ws.Range("C20").Value = 0.0172824121636026 ; $ws.Range("D20").Value = 0.2028311534361364 ; $ws.Range("E20").Value = 0.3393607109364183 ; $ws.Range("F20").Value = 1.013598427796182 ; $ws.Range("G20").Value = 0.8047071325291313 ; $ws.Range("H20").Value = 0.8068796609075264 ; $ws.Range("I20").Value = 0.7169349375235328 ; $ws.Range("J20").Value = 0.5736832052066347 ; $ws.Range("M20").Value = 15.55090899995321
$ws.Range("C21").Value = 0.01822742383949105 ; $ws.Range("D21").Value = 0.2173962949703707 ; $ws.Range("E21").Value = 0.3731651796871915 ; $ws.Range("F21").Value = 0.9665771928908526 ; $ws.Range("G21").Value = 0.7993297783727655 ; $ws.Range("H21").Value = 0.7861122440234567 ; $ws.Range("I21").Value = 0.7116480480232639 ; $ws.Range("J21").Value = 0.6371341232034808 ; $ws.Range("M21").Value = 17.38971439163328
$ws.Range("C22").Value = 0.01886327734367654 ; $ws.Range("D22").Value = 0.2269361570177182 ; $ws.Range("E22").Value = 0.3954411054839397 ; $ws.Range("F22").Value = 0.9388962683769932 ; $ws.Range("G22").Value = 0.798417994529899 ; $ws.Range("H22").Value = 0.7741278856351812 ; $ws.Range("I22").Value = 0.7100200064232354 ; $ws.Range("J22").Value = 0.6790219651729785 ; $ws.Range("M22").Value = 18.59156285596544
$ws.Range("C23").Value = 0.01852218386840576 ; $ws.Range("D23").Value = 0.221842590145144 ; $ws.Range("E23").Value = 0.3835344247531225 ; $ws.Range("F23").Value = 0.9533850074419234 ; $ws.Range("G23").Value = 0.798660564670655 ; $ws.Range("H23").Value = 0.7803761612723008 ; $ws.Range("I23").Value = 0.7107179737914464 ; $ws.Range("J23").Value = 0.6566250831826892 ; $ws.Range("M23").Value = 17.95007797062914
$ws.Range("C24").Value = 0.01726726164812931 ; $ws.Range("D24").Value = 0.2025934290217037 ; $ws.Range("E24").Value = 0.3388110898253274 ; $ws.Range("F24").Value = 1.014417137021965 ; $ws.Range("G24").Value = 0.8048376549386376 ; $ws.Range("H24").Value = 0.8072453202719316 ; $ws.Range("I24").Value = 0.7170513690931273 ; $ws.Range("J24").Value = 0.5726527261905687 ; $ws.Range("M24").Value = 15.52085000930197
$ws.Range("C25").Value = 0.01598637786256774 ; $ws.Range("D25").Value = 0.1819130932251483 ; $ws.Range("E25").Value = 0.2912714730620536 ; $ws.Range("F25").Value = 1.092756530483513 ; $ws.Range("G25").Value = 0.8220925342952512 ; $ws.Range("H25").Value = 0.8427851921260441 ; $ws.Range("I25").Value = 0.7313588421243438 ; $ws.Range("J25").Value = 0.4836633939982278 ; $ws.Range("M25").Value = 12.89824806121686
